# Weekly update: insert a new price record as the new first data row (142)
# for this Pepino ensalada subset, pushing all existing rows (142-169) down
# by one (they end up at 143-170). The new row carries the same constant
# attributes (market, region, category, etc.) as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 142; Excel shifts rows 142:169 down to 143:170
# and keeps the dimension/formatting (e.g. the date style on column D) of
# the row above, matching the existing data rows.
$ws.Rows("142:142").Insert()

$ws.Range("A142").Value = 7
$ws.Range("B142").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C142").Value = "Ñuble"
$ws.Range("D142").Value = 44508
$ws.Range("E142").Value = 16
$ws.Range("F142").Value = 100112043
$ws.Range("G142").Value = "Pepino ensalada"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 100
$ws.Range("K142").Value = 9000
$ws.Range("L142").Value = 10000
$ws.Range("M142").Value = 9500
$ws.Range("N142").Value = "$/caja 80 unidades"
$ws.Range("O142").Value = "Región del Maule"
$ws.Range("P142").Value = 119
$ws.Range("Q142").Value = 80
$ws.Range("R142").Value = "Hortaliza"
